# Regenerate merged AHB files
# 1. Rename the comparison-column headers:
#      "<Name>_old" -> "<Name>_FV2310"
#      "<Name>_new" -> "<Name>_FV2404"
# 2. Turn the data range into an Excel Table ("Table1").
# 3. Freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J (1-10): "_old" suffix becomes "_FV2310"
for ($i = 0; $i -lt $baseHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseHeaders[$i] + "_FV2310"
}

# Column K (11) stays "diff" - untouched.

# Columns L-U (12-21): "_new" suffix becomes "_FV2404"
for ($i = 0; $i -lt $baseHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $baseHeaders[$i] + "_FV2404"
}

# Convert the used range into a table, matching the dimension A1:U85.
$tableRange = $ws.Range("A1:U85")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze the top (header) row.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
